$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as text so values like "1.00" or "5.31" are not
# silently converted into numbers by Excel when assigned via .Value
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "54.285.72"
$ws.Range("E2").Value = "  -6.59%  "

# Row 3
$ws.Range("D3").Value = "2.455.01"
$ws.Range("E3").Value = "  -8.49%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.19%  "

# Row 5
$ws.Range("D5").Value = "465.88"
$ws.Range("E5").Value = "  -5.73%  "

# Row 6
$ws.Range("D6").Value = "132.67"
$ws.Range("E6").Value = "  -2.87%  "

# Row 7
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.13%  "

# Row 8
$ws.Range("D8").Value = "0.489"
$ws.Range("E8").Value = "  -6.25%  "

# Row 9
$ws.Range("D9").Value = "2.455.91"
$ws.Range("E9").Value = "  -8.68%  "

# Row 10
$ws.Range("D10").Value = "0.0953"
$ws.Range("E10").Value = "  -5.09%  "

# Row 11
$ws.Range("D11").Value = "5.31"
$ws.Range("E11").Value = "  -9.13%  "

# Row 12
$ws.Range("D12").Value = "0.317"
$ws.Range("E12").Value = "  -5.47%  "

# Row 13
$ws.Range("E13").Value = "  -3.61%  "

# Row 14
$ws.Range("D14").Value = "2.884.83"
$ws.Range("E14").Value = "  -8.91%  "

# Row 15
$ws.Range("D15").Value = "54.378.52"
$ws.Range("E15").Value = "  -6.36%  "

# Row 16
$ws.Range("D16").Value = "19.80"
$ws.Range("E16").Value = "  -5.01%  "

# Row 17
$ws.Range("D17").Value = "0.0000131"
$ws.Range("E17").Value = "  -0.51%  "

# Row 18
$ws.Range("D18").Value = "2.465.43"
$ws.Range("E18").Value = "  -8.02%  "

# Row 19
$ws.Range("D19").Value = "4.20"
$ws.Range("E19").Value = "  -8.36%  "

# Row 20
$ws.Range("D20").Value = "307.15"
$ws.Range("E20").Value = "  -7.96%  "

# Row 21
$ws.Range("D21").Value = "9.45"
$ws.Range("E21").Value = "  -11.12%  "

# Row 22
$ws.Range("D22").Value = "0.999"

# Row 23
$ws.Range("D23").Value = "5.67"
$ws.Range("E23").Value = "  +1.20%  "

# Row 24
$ws.Range("D24").Value = "5.35"
$ws.Range("E24").Value = "  -11.72%  "

# Row 25
$ws.Range("D25").Value = "56.66"
$ws.Range("E25").Value = "  -7.88%  "

# Row 26
$ws.Range("E26").Value = "  +0.81%  "

# Row 27
$ws.Range("D27").Value = "0.385"
$ws.Range("E27").Value = "  -6.87%  "

# Row 28
$ws.Range("D28").Value = "2.556.04"
$ws.Range("E28").Value = "  -9.06%  "

# Row 29
$ws.Range("D29").Value = "0.153"
$ws.Range("E29").Value = "  -10.07%  "

# Row 30
$ws.Range("D30").Value = "7.22"
$ws.Range("E30").Value = "  -0.16%  "

# Row 31
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.16%  "

# Row 32
$ws.Range("D32").Value = "0.0₃0717"
$ws.Range("E32").Value = "  -10.87%  "

# Row 33
$ws.Range("D33").Value = "147.01"
$ws.Range("E33").Value = "  +0.10%  "

# Row 34
$ws.Range("D34").Value = "17.80"
$ws.Range("E34").Value = "  -4.26%  "

# Row 35
$ws.Range("D35").Value = "1.43"
$ws.Range("E35").Value = "  -8.46%  "

# Row 36
$ws.Range("D36").Value = "5.03"
$ws.Range("E36").Value = "  -2.77%  "

# Row 37
$ws.Range("D37").Value = "3.55"
$ws.Range("E37").Value = "  -12.38%  "

# Row 38
$ws.Range("E38").Value = "  -3.09%  "

# Row 39
$ws.Range("D39").Value = "0.800"
$ws.Range("E39").Value = "  -10.56%  "

# Row 40
$ws.Range("D40").Value = "33.61"
$ws.Range("E40").Value = "  -6.26%  "

# Row 41
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.30%  "

# Row 42
$ws.Range("D42").Value = "0.601"
$ws.Range("E42").Value = "  +2.32%  "

# Row 43
$ws.Range("D43").Value = "0.0527"
$ws.Range("E43").Value = "  -2.55%  "

# Row 44
$ws.Range("D44").Value = "3.28"
$ws.Range("E44").Value = "  -3.88%  "

# Row 45
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").Value = "10.22"
$ws.Range("E45").Value = "  -1.15%  "

# Row 46
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "1.24"
$ws.Range("E46").Value = "  -7.41%  "

# Row 47
$ws.Range("D47").Value = "1.954.15"
$ws.Range("E47").Value = "  -8.12%  "

# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "0.0874"
$ws.Range("E48").Value = "  +0.49%  "

# Row 49
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0218"
$ws.Range("E49").Value = "  -0.58%  "

# Row 50
$ws.Range("D50").Value = "4.22"
$ws.Range("E50").Value = "  -7.01%  "

# Row 51
$ws.Range("D51").Value = "16.61"
$ws.Range("E51").Value = "  -9.02%  "

